$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last populated row and append the new data row right after it
# (the sheet is a simple flat table: date, volume, high, low, open,
# close, adj_close, ticker in columns A:H).
$usedRange = $ws.UsedRange
$srcRow = $usedRange.Row + $usedRange.Rows.Count - 1
$newRow = $srcRow + 1

# Copy the formatting of the row above down onto the new row first
# (so the date cell picks up the same date number format / style as
# the rest of column A).
$ws.Range("A$srcRow`:H$srcRow").Copy()
$ws.Range("A$newRow`:H$newRow").PasteSpecial(-4122)  # xlPasteFormats

# The new row's adj_close/ticker text ("0.709999978542328" / "BWZ.MI")
# happen to match the row above exactly, so copy those two cells'
# values straight down - this reuses the existing shared-string
# entries for those text cells without perturbing the style table.
$ws.Range("G$srcRow`:H$srcRow").Copy()
$ws.Range("G$newRow`:H$newRow").PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = 0

# Fill in the row-specific numeric data (date serial, volume, high,
# low, open, close).
$ws.Cells.Item($newRow, 1).Value2 = 45449.2916666667
$ws.Cells.Item($newRow, 2).Value2 = 2455
$ws.Cells.Item($newRow, 3).Value2 = 0.730000019073486
$ws.Cells.Item($newRow, 4).Value2 = 0.709999978542328
$ws.Cells.Item($newRow, 5).Value2 = 0.714999973773956
$ws.Cells.Item($newRow, 6).Value2 = 0.709999978542328
